$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-22 05:48:39'
$ws.Range("E3").Value = '2026-02-22 05:48:42'
$ws.Range("L3").Value = '23.4 km/h - 90º 5:12 TU'
$ws.Range("O3").Value = '2.9 °C'
$ws.Range("E4").Value = '2026-02-22 05:48:45'
$ws.Range("H4").Value = '''80%'
$ws.Range("N4").Value = '3.5 °C 5:18 TU'
$ws.Range("O4").Value = '6.3 °C'
$ws.Range("E5").Value = '2026-02-22 05:48:48'
$ws.Range("L5").Value = '7.9 km/h - 236º 5:11 TU'
$ws.Range("E6").Value = '2026-02-22 05:48:50'
$ws.Range("E7").Value = '2026-02-22 05:48:53'
$ws.Range("H7").Value = '''63%'
$ws.Range("N7").Value = '10.7 °C 5:19 TU'
$ws.Range("O7").Value = '11.8 °C'
$ws.Range("E8").Value = '2026-02-22 05:48:56'
$ws.Range("H8").Value = '''51%'
$ws.Range("J8").Value = '1027.7 hPa'
$ws.Range("M8").Value = '14.5 °C 5:29 TU'
$ws.Range("O8").Value = '12.2 °C'
$ws.Range("E9").Value = '2026-02-22 05:48:59'
$ws.Range("O9").Value = '4.4 °C'
$ws.Range("E10").Value = '2026-02-22 05:49:01'
$ws.Range("N10").Value = '2.9 °C 5:29 TU'
$ws.Range("E11").Value = '2026-02-22 05:49:04'
$ws.Range("N11").Value = '0.9 °C 5:22 TU'
$ws.Range("O11").Value = '1.5 °C'
$ws.Range("E12").Value = '2026-02-22 05:49:07'
$ws.Range("O12").Value = '4.5 °C'
$ws.Range("E13").Value = '2026-02-22 05:49:09'
$ws.Range("N13").Value = '-3.5 °C 5:01 TU'
$ws.Range("O13").Value = '-2.5 °C'
$ws.Range("E14").Value = '2026-02-22 05:49:12'
$ws.Range("N14").Value = '6.7 °C 5:29 TU'
$ws.Range("O14").Value = '7.5 °C'
$ws.Range("E15").Value = '2026-02-22 05:49:15'
$ws.Range("H15").Value = '''87%'
$ws.Range("O15").Value = '4.4 °C'
$ws.Range("E16").Value = '2026-02-22 05:49:17'
$ws.Range("H16").Value = '''22%'
$ws.Range("L16").Value = '15.5 km/h - 214º 5:04 TU'
$ws.Range("M16").Value = '5.6 °C 5:23 TU'
$ws.Range("O16").Value = '4.0 °C'
$ws.Range("E17").Value = '2026-02-22 05:49:20'
$ws.Range("M17").Value = '8.5 °C 5:07 TU'
$ws.Range("E18").Value = '2026-02-22 05:49:23'
$ws.Range("N18").Value = '0.1 °C 5:25 TU'
$ws.Range("O18").Value = '1.3 °C'
$ws.Range("E19").Value = '2026-02-22 05:49:25'
$ws.Range("N19").Value = '6.4 °C 5:00 TU'
$ws.Range("O19").Value = '7.8 °C'
$ws.Range("E20").Value = '2026-02-22 05:49:28'
$ws.Range("O20").Value = '0.4 °C'
$ws.Range("E21").Value = '2026-02-22 05:49:31'
$ws.Range("H21").Value = '''78%'
$ws.Range("J21").Value = '1033.5 hPa'
$ws.Range("N21").Value = '0.8 °C 5:28 TU'
$ws.Range("O21").Value = '2.5 °C'
$ws.Range("E22").Value = '2026-02-22 05:49:34'
$ws.Range("H22").Value = '''25%'
$ws.Range("M22").Value = '4.3 °C 5:19 TU'
$ws.Range("O22").Value = '3.2 °C'
$ws.Range("E23").Value = '2026-02-22 05:49:36'
$ws.Range("E24").Value = '2026-02-22 05:49:39'
$ws.Range("J24").Value = '1031.4 hPa'
$ws.Range("E25").Value = '2026-02-22 05:49:42'
$ws.Range("H25").Value = '''28%'
$ws.Range("E26").Value = '2026-02-22 05:49:44'
$ws.Range("E27").Value = '2026-02-22 05:49:47'
$ws.Range("H27").Value = '''30%'
$ws.Range("L27").Value = '13.7 km/h - 230º 5:08 TU'
$ws.Range("M27").Value = '5.2 °C 5:10 TU'
$ws.Range("O27").Value = '4.0 °C'
$ws.Range("E28").Value = '2026-02-22 05:49:50'
$ws.Range("O28").Value = '2.0 °C'
$ws.Range("E29").Value = '2026-02-22 05:49:53'
$ws.Range("N29").Value = '2.5 °C 5:26 TU'
$ws.Range("O29").Value = '4.4 °C'
$ws.Range("E30").Value = '2026-02-22 05:49:55'
$ws.Range("N30").Value = '6.0 °C 5:12 TU'
$ws.Range("O30").Value = '7.7 °C'
$ws.Range("E31").Value = '2026-02-22 05:49:58'
$ws.Range("H31").Value = '''65%'
$ws.Range("N31").Value = '10.1 °C 5:16 TU'
$ws.Range("O31").Value = '11.9 °C'
$ws.Range("E32").Value = '2026-02-22 05:50:01'
$ws.Range("K32").Value = '-0.1 MJ/m2'
$ws.Range("E33").Value = '2026-02-22 05:50:03'
$ws.Range("N33").Value = '-0.3 °C 5:29 TU'
$ws.Range("O33").Value = '1.2 °C'
$ws.Range("E34").Value = '2026-02-22 05:50:06'
$ws.Range("H34").Value = '''44%'
$ws.Range("O34").Value = '2.5 °C'
$ws.Range("E35").Value = '2026-02-22 05:50:08'
$ws.Range("H35").Value = '''40%'
$ws.Range("K35").Value = '-0.1 MJ/m2'
$ws.Range("O35").Value = '5.7 °C'
$ws.Range("E36").Value = '2026-02-22 05:50:10'
$ws.Range("E37").Value = '2026-02-22 05:50:12'
$ws.Range("O37").Value = '-0.5 °C'
$ws.Range("E38").Value = '2026-02-22 05:50:15'
$ws.Range("E39").Value = '2026-02-22 05:50:18'
$ws.Range("G39").Value = '53 cm'
$ws.Range("H39").Value = '''27%'
$ws.Range("I39").Value = '0.0 mm'
$ws.Range("K39").Value = '-0.1 MJ/m2'
$ws.Range("L39").Value = '16.9 km/h - 247º 3:58 TU'
$ws.Range("M39").Value = '4.8 °C 0:28 TU'
$ws.Range("N39").Value = '3.1 °C 4:50 TU'
$ws.Range("O39").Value = '4.0 °C'
$ws.Range("E40").Value = '2026-02-22 05:50:21'
$ws.Range("E41").Value = '2026-02-22 05:50:23'
$ws.Range("J41").Value = '1028.4 hPa'
$ws.Range("O41").Value = '5.2 °C'
$ws.Range("E42").Value = '2026-02-22 05:50:26'
$ws.Range("N42").Value = '3.1 °C 5:13 TU'
$ws.Range("O42").Value = '4.8 °C'
$ws.Range("E43").Value = '2026-02-22 05:50:29'
$ws.Range("N43").Value = '0.4 °C 5:21 TU'
$ws.Range("O43").Value = '2.0 °C'
$ws.Range("E44").Value = '2026-02-22 05:50:31'
$ws.Range("N44").Value = '-1.3 °C 5:25 TU'
$ws.Range("O44").Value = '0.0 °C'
$ws.Range("E45").Value = '2026-02-22 05:50:34'
$ws.Range("J45").Value = '1031.6 hPa'
$ws.Range("N45").Value = '2.0 °C 5:28 TU'
$ws.Range("O45").Value = '4.5 °C'
$ws.Range("E46").Value = '2026-02-22 05:50:37'
$ws.Range("N46").Value = '0.0 °C 5:29 TU'
$ws.Range("O46").Value = '1.2 °C'
